# AB#44 — fixed minor export issues related to data reloading.
#
# The "Maximum force [N]" measurement column in the TestSubject data-export
# template is replaced by a combined "Maximum contraction [Nm] / Fatigue [%]"
# column. Renaming the table header cell automatically renames the bound
# Excel Table column (MeasurementsTable) and rewrites the shared-string
# table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSubject")

# Rename the 5th measurement column header (currently "Maximum force [N]")
# to the new combined contraction/fatigue header. The embedded line break
# is required so Excel wraps the header onto two lines, matching the
# template's tall header row.
$ws.Range("E6").Value = "Maximum contraction [Nm] " + [char]10 + "/ Fatigue [%] "

# Widen column E so the new, longer two-line header remains readable.
$ws.Columns("E").ColumnWidth = 28.3

# Restore the cursor/selection that was active when the template was last
# saved.
[void]$ws.Range("G15").Select()
